$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.989.54"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "'2.947.44"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'379.78"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'101.00"
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").Value = "'36.16"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'3.403.57"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'18.30"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").Value = "'11.94"
$ws.Range("E16").Value = "  +68.27%  "
$ws.Range("D17").Value = "'2.958.82"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "'50.964.51"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("E20").Value = "  -4.71%  "
$ws.Range("D21").Value = "'12.41"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D23").Value = "'69.50"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'266.85"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "'3.21"
$ws.Range("E25").Value = "  +12.04%  "
$ws.Range("D26").Value = "'8.19"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'7.15"
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'25.61"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("D32").Value = "'10.06"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").Value = "'50.50"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'33.53"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -5.35%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'3.09"
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'2.52"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.81"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").Value = "'119.72"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").Value = "'21.45"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "'3.48"
$ws.Range("E45").Value = "  +7.74%  "
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D48").Value = "'2.009.99"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'0.262"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("E50").Value = "  -6.12%  "
$ws.Range("D51").Value = "'5.29"
$ws.Range("E51").Value = "  +4.46%  "
